$d = $word.ActiveDocument

# --- Paragraph 1: pPr changes (border + indent) ---
$p1 = $d.Paragraphs(1)
$pf1 = $p1.Range.ParagraphFormat

# w:ind w:left="120" -> w:left="225" (OOXML twips -> COM points: value/20)
$pf1.LeftIndent = 225 / 20

# Add a paragraph border (top/left/bottom/right), 5-twip space, no line -
# matches <w:pBdr><w:top w:space="5"/><w:left .../><w:bottom .../><w:right .../></w:pBdr>
$borders1 = $pf1.Borders
$borders1.DistanceFromTop = 5
$borders1.DistanceFromBottom = 5
$borders1.DistanceFromLeft = 5
$borders1.DistanceFromRight = 5

# --- Paragraph 1: replace the placeholder ID text ---
$rng1 = $p1.Range
[void]$rng1.Find.Execute("**ID__AFFARS_pgi_5308_topic_3__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5308_404_90__ID**", 2)

# --- Paragraph 1: drop the trailing " " run that used to follow the ID text ---
$p1b = $d.Paragraphs(1)
$text1b = $p1b.Range.Text
# $text1b ends with the paragraph mark; the character before it is the space.
$spaceStart = $p1b.Range.Start + $text1b.Length - 2
$spaceEnd = $spaceStart + 1
$spaceRange = $d.Range($spaceStart, $spaceEnd)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}
